$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set marketDaysMode (B5) from "Auto" to "Manual"
$ws.Range("B5").Value = "Manual"

# Set manualMarketDays (B6) from 365 to 250
$ws.Range("B6").Value = 250

# Update the selected cell on the sheet to B7 (as recorded in the saved file)
$ws.Range("B7").Select()
